$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New experiment data (date 9.25.2021) appended as rows 172-226.
# Columns: A=date, B=strain, C=media, D=strainmedia (=CONCATENATE(B,C)),
#          E=replicate, F=fitness
# ---------------------------------------------------------------------------
$data = @(
    @{row=172; strain="ancestor"; media="all"; rep=2; fit=0.04}
    @{row=173; strain="ancestor"; media="all"; rep=3; fit=-0.01}
    @{row=174; strain="ancestor"; media="glucose"; rep=1; fit=-0.05}
    @{row=175; strain="ancestor"; media="glucose"; rep=2; fit=0.09}
    @{row=176; strain="ancestor"; media="glucose"; rep=3; fit=0.11}
    @{row=177; strain="ancestor"; media="aminoacids"; rep=1; fit=-0.09}
    @{row=178; strain="ancestor"; media="aminoacids"; rep=2; fit=0.13}
    @{row=179; strain="ancestor"; media="aminoacids"; rep=3; fit=0.13}
    @{row=180; strain="ancestor"; media="lactate"; rep=1; fit=0.02}
    @{row=181; strain="ancestor"; media="lactate"; rep=2; fit=-0.12}
    @{row=182; strain="ancestor"; media="lactate"; rep=3; fit=-0.11}
    @{row=183; strain="morA"; media="all"; rep=1; fit=0.01}
    @{row=184; strain="morA"; media="all"; rep=2; fit=0.03}
    @{row=185; strain="morA"; media="all"; rep=3; fit=0.04}
    @{row=186; strain="morA"; media="glucose"; rep=1; fit=0.12}
    @{row=187; strain="morA"; media="glucose"; rep=2; fit=0.14000000000000001}
    @{row=188; strain="morA"; media="glucose"; rep=3; fit=0.02}
    @{row=189; strain="morA"; media="aminoacids"; rep=1; fit=-0.36}
    @{row=190; strain="morA"; media="aminoacids"; rep=2; fit=-0.61}
    @{row=191; strain="morA"; media="aminoacids"; rep=3; fit=-0.43}
    @{row=192; strain="morA"; media="lactate"; rep=1; fit=0.21}
    @{row=193; strain="morA"; media="lactate"; rep=2; fit=-0.55000000000000004}
    @{row=194; strain="morA"; media="lactate"; rep=3; fit=0.18}
    @{row=195; strain="morA"; media="half all"; rep=1; fit=0}
    @{row=196; strain="morA"; media="half all"; rep=2; fit=-0.13}
    @{row=197; strain="morA"; media="half all"; rep=3; fit=-0.09}
    @{row=198; strain="morA"; media="half glucose"; rep=1; fit=0.14000000000000001}
    @{row=199; strain="morA"; media="half glucose"; rep=2; fit=0.14000000000000001}
    @{row=200; strain="morA"; media="half glucose"; rep=3; fit=0.32}
    @{row=201; strain="morA"; media="double aminoacids"; rep=1; fit=-0.11}
    @{row=202; strain="morA"; media="double aminoacids"; rep=2; fit=-0.22}
    @{row=203; strain="morA"; media="double aminoacids"; rep=3; fit=0.02}
    @{row=204; strain="lasR R216Q"; media="all"; rep=1; fit=0.18}
    @{row=205; strain="lasR R216Q"; media="all"; rep=2; fit=0.37}
    @{row=206; strain="lasR R216Q"; media="all"; rep=3; fit=0.11}
    @{row=207; strain="lasR R216Q"; media="glucose"; rep=1; fit=0.41}
    @{row=208; strain="lasR R216Q"; media="glucose"; rep=2; fit=0.34}
    @{row=209; strain="lasR R216Q"; media="glucose"; rep=3; fit=0.19}
    @{row=210; strain="lasR R216Q"; media="aminoacids"; rep=1; fit=0.64}
    @{row=211; strain="lasR R216Q"; media="aminoacids"; rep=2; fit=0.56000000000000005}
    @{row=212; strain="lasR R216Q"; media="aminoacids"; rep=3; fit=0.45}
    @{row=213; strain="lasR R216Q"; media="lactate"; rep=1; fit=0.63}
    @{row=214; strain="lasR R216Q"; media="lactate"; rep=3; fit=0.37}
    @{row=215; strain="ΔPA14_45920..PA14_46440"; media="all"; rep=1; fit=0.38}
    @{row=216; strain="ΔPA14_45920..PA14_46440"; media="all"; rep=2; fit=0.24}
    @{row=217; strain="ΔPA14_45920..PA14_46440"; media="all"; rep=3; fit=0.14000000000000001}
    @{row=218; strain="ΔPA14_45920..PA14_46440"; media="glucose"; rep=1; fit=0.17}
    @{row=219; strain="ΔPA14_45920..PA14_46440"; media="glucose"; rep=2; fit=0.13}
    @{row=220; strain="ΔPA14_45920..PA14_46440"; media="glucose"; rep=3; fit=0.15}
    @{row=221; strain="ΔPA14_45920..PA14_46440"; media="aminoacids"; rep=1; fit=0.28999999999999998}
    @{row=222; strain="ΔPA14_45920..PA14_46440"; media="aminoacids"; rep=2; fit=0.19}
    @{row=223; strain="ΔPA14_45920..PA14_46440"; media="aminoacids"; rep=3; fit=0.16}
    @{row=224; strain="ΔPA14_45920..PA14_46440"; media="lactate"; rep=1; fit=0.61}
    @{row=225; strain="ΔPA14_45920..PA14_46440"; media="lactate"; rep=2; fit=0.25}
    @{row=226; strain="ΔPA14_45920..PA14_46440"; media="lactate"; rep=3; fit=0.82}
)

$firstRow = 172
$lastRow = 226

# Copy the formatting (font/alignment/number-format) of the last "real" data
# row (130) down across the newly-populated rows, then copy that same row's
# column-A formatting so the new text date keeps the normal "all" style
# (rather than being auto-parsed as a serial date number).
$ws.Range("A130:F130").Copy() | Out-Null
$ws.Range("A$firstRow`:F$lastRow").PasteSpecial(-4122) | Out-Null

foreach ($row in $data) {
    $r = $row.row

    # Force column A to stay text ("9.25.2021") instead of being reinterpreted
    # as a date serial number.
    $ws.Cells.Item($r, 1).Value = "'9.25.2021"

    $ws.Cells.Item($r, 2).Value = $row.strain
    $ws.Cells.Item($r, 3).Value = $row.media
    $ws.Cells.Item($r, 4).Formula = "=CONCATENATE(B$r,C$r)"
    $ws.Cells.Item($r, 5).Value = $row.rep
    $ws.Cells.Item($r, 6).Value = $row.fit
}

# Re-apply row 130's formatting to column A only, restoring the shared style
# (s="2") that the quoted-text assignment above replaced with a one-off
# quote-prefixed style.
$ws.Range("A130").Copy() | Out-Null
$ws.Range("A$firstRow`:A$lastRow").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Trailing blank rows (227-241), matching the old trailing blank rows that
# used to sit right after the data (previously rows 172-179, now pushed down
# to make room for the new data block above).
# ---------------------------------------------------------------------------
$ws.Range("E130").Copy() | Out-Null
$ws.Range("E227:E241").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Window / selection state, matching the author re-scrolling down to the
# newly entered data and leaving the cursor on D216.
# ---------------------------------------------------------------------------
$ws.Range("D216").Select() | Out-Null
